# Insert a new weekly record at row 123 ("Hortaliza, Macroferia Regional de
# Talca - Acelga"): push the existing rows 123-221 down to 124-222 and
# populate the new row 123 with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 123..221 down to 124..222, opening up row 123.
$ws.Rows.Item(123).EntireRow.Insert()

# Fill in the newly opened row 123 with the new weekly data point.
$ws.Cells.Item(123, 1).Value  = 5
$ws.Cells.Item(123, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(123, 3).Value  = 'Maule'
$ws.Cells.Item(123, 4).Value  = 44651
$ws.Cells.Item(123, 5).Value  = 7
$ws.Cells.Item(123, 6).Value  = 100112009
$ws.Cells.Item(123, 7).Value  = 'Acelga'
$ws.Cells.Item(123, 8).Value  = 'Sin especificar'
$ws.Cells.Item(123, 9).Value  = 'Primera'
$ws.Cells.Item(123, 10).Value = 350
$ws.Cells.Item(123, 11).Value = 3500
$ws.Cells.Item(123, 12).Value = 3500
$ws.Cells.Item(123, 13).Value = 3500
$ws.Cells.Item(123, 14).Value = '$/docena de atados (4 kilos)'
$ws.Cells.Item(123, 15).Value = 'Región del Maule'
$ws.Cells.Item(123, 16).Value = 875
$ws.Cells.Item(123, 17).Value = 4
$ws.Cells.Item(123, 18).Value = 'Hortaliza'
